$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column D (rows 3-18, skipping row 12) to use formulas "=8+<old value>"
$ws.Range("D3").Formula  = "=8+6"
$ws.Range("D4").Formula  = "=8+8"
$ws.Range("D5").Formula  = "=8+12"
$ws.Range("D6").Formula  = "=8+14"
$ws.Range("D7").Formula  = "=8+16"
$ws.Range("D8").Formula  = "=8+18"
$ws.Range("D9").Formula  = "=8+20"
$ws.Range("D10").Formula = "=8+22"
$ws.Range("D11").Formula = "=8+24"
# D12 unchanged (stays 26)
$ws.Range("D13").Formula = "=8+28"
$ws.Range("D14").Formula = "=8+30"
$ws.Range("D15").Formula = "=8+32"
$ws.Range("D16").Formula = "=8+4"
$ws.Range("D17").Formula = "=8+4"
$ws.Range("D18").Formula = "=8+4"

# Add new instruction-memory rows describing the factorial accelerator routine
# (values entered in this order so new shared-string indices line up with the
# original authoring session)
$ws.Range("D22").Value = "load n"
$ws.Range("D23").Value = "load go"
$ws.Range("D24").Value = "send n"
$ws.Range("D25").Value = "send go"
$ws.Range("D27").Value = "read done"
$ws.Range("D30").Value = "read nf"
$ws.Range("D31").Value = "store nf"
$ws.Range("D28").Value = "branch if done or error"

$ws.Range("C22").Value = "load n"
$ws.Range("C23").Value = "branch to factorial"
$ws.Range("C27").Value = "wait till return"
$ws.Range("C30").Value = "store return value"

# Update the sheet view to match what was recorded during editing
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("I28").Select()
